# REVER_DailyTracker - add new rows to the OCT-2020 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OCT-2020")

# Row 20 is the "Week off" template row (A:1 B:date C:- D:"Week off" E:- F:- G:-)
# Copy it down to rows 26, 27 and 28 so formatting/styles match exactly.
$ws.Range("A20:G20").Copy($ws.Range("A26:G26"))
$ws.Range("A20:G20").Copy($ws.Range("A27:G27"))
$ws.Range("A20:G20").Copy($ws.Range("A28:G28"))

# Row 24 is the "QMVAR / issue fixing / WIP" template row
# Copy it down to rows 29, 30 and 31.
$ws.Range("A24:G24").Copy($ws.Range("A29:G29"))
$ws.Range("A24:G24").Copy($ws.Range("A30:G30"))
$ws.Range("A24:G24").Copy($ws.Range("A31:G31"))

# Fill in the per-row numbers / dates
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 44128

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 44129

$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "10/126/2020"

$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 44126
$ws.Cells.Item(29, 4).Value = "issue fixing"
$ws.Cells.Item(29, 5).ClearContents()

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 44126
$ws.Cells.Item(30, 4).Value = "issue fixing"
$ws.Cells.Item(30, 5).ClearContents()

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 44126
$ws.Cells.Item(31, 4).Value = "issue fixing"
$ws.Cells.Item(31, 5).ClearContents()

# Update the visible selection / active cell like the author left it
$ws.Range("D33").Select()
